$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.089.29"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "1.835.51"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.04"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6146"
$ws.Range("E6").Value = "  -2.66%  "
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07466"
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2919"
$ws.Range("E9").Value = "  -0.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "23.10"
$ws.Range("E10").Value = "  +0.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07685"
$ws.Range("E11").Value = "  -0.31%  "
$ws.Range("D12").Value = "1.829.22"
$ws.Range("E12").Value = "  -0.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.002"
$ws.Range("E13").Value = "  +0.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6713"
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "82.60"
$ws.Range("E15").Value = "  -0.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000009215"
$ws.Range("E16").Value = "  -4.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.923"
$ws.Range("E17").Value = "  -2.51%  "
$ws.Range("D18").Value = "29.059.51"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "2.087.04"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "231.21"
$ws.Range("E20").Value = "  +2.07%  "
$ws.Range("E21").Value = "  +0.69%  "
$ws.Range("E22").Value = "  +0.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.183"
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("E24").Value = "  +0.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "159.67"
$ws.Range("E25").Value = "  -0.25%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1386"
$ws.Range("E26").Value = "  -1.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.494"
$ws.Range("E27").Value = "  -0.50%  "
$ws.Range("E28").Value = "  -0.64%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.156"
$ws.Range("E30").Value = "  +0.86%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.133"
$ws.Range("E31").Value = "  +1.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.05514"
$ws.Range("E32").Value = "  +2.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.209"
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7441"
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.835"
$ws.Range("E35").Value = "  -1.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.139"
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.661"
$ws.Range("E37").Value = "  +0.28%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.776"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01779"
$ws.Range("E39").Value = "  -0.55%  "
$ws.Range("D40").Value = "1.213.38"
$ws.Range("E40").Value = "  -2.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.476"
$ws.Range("E41").Value = "  -2.38%  "
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "102.04"
$ws.Range("E44").Value = "  +0.37%  "
$ws.Range("D45").Value = "1.985.73"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "65.56"
$ws.Range("E46").Value = "  +1.17%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000123"
$ws.Range("E47").Value = "  -0.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5093"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4063"
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.099"
$ws.Range("E50").Value = "  +1.20%  "
$ws.Range("E51").Value = "  +1.08%  "
